$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 403.8125
$ws.Range("I28").Value = 403.8125
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 403.8125
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 81.1875
$ws.Range("N28").Value = $null

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 9500
$ws.Range("I47").Value = 9500
$ws.Range("K47").Value = 9500
$ws.Range("M47").Value = -8528

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 402.9
$ws.Range("I107").Value = 381.55554
$ws.Range("J107").Value = 595
$ws.Range("K107").Value = 381.55554
$ws.Range("L107").Value = 595
$ws.Range("M107").Value = 1538.44446
$ws.Range("N107").Value = -4435

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2706
$ws.Range("I125").Value = 700
$ws.Range("J125").Value = 3207.5
$ws.Range("K125").Value = 6300
$ws.Range("L125").Value = 28867.5
$ws.Range("M125").Value = -3840
$ws.Range("N125").Value = -33787.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1084.3846
$ws.Range("I127").Value = 497
$ws.Range("J127").Value = 1133.3334
$ws.Range("K127").Value = 1491
$ws.Range("L127").Value = 3400.0002
$ws.Range("M127").Value = 3469
$ws.Range("N127").Value = -13320.0002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5702.5
$ws.Range("I132").Value = 4909.9395
$ws.Range("J132").Value = 7079.0527
$ws.Range("K132").Value = 14729.8185
$ws.Range("L132").Value = 21237.1581
$ws.Range("M132").Value = -12199.8185
$ws.Range("N132").Value = -26297.1581

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1909.2632
$ws.Range("I137").Value = 1879.7037
$ws.Range("J137").Value = 1981.8182
$ws.Range("K137").Value = 5639.1111
$ws.Range("L137").Value = 5945.4546
$ws.Range("M137").Value = -3089.1111
$ws.Range("N137").Value = -11045.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 50000
$ws.Range("J7").Value = 50000
$ws.Range("L7").Value = 50000
$ws.Range("N7").Value = -50228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7736.892
$ws.Range("I32").Value = 8269.603999999999
$ws.Range("J32").Value = 5805.8125
$ws.Range("K32").Value = 8269.603999999999
$ws.Range("L32").Value = 5805.8125
$ws.Range("M32").Value = -7982.603999999999
$ws.Range("N32").Value = -6379.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1239.75
$ws.Range("I110").Value = 952.2692
$ws.Range("J110").Value = 2485.5
$ws.Range("K110").Value = 952.2692
$ws.Range("L110").Value = 2485.5
$ws.Range("M110").Value = 1092.7308
$ws.Range("N110").Value = -6575.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 33000
$ws.Range("J121").Value = 33000
$ws.Range("L121").Value = 33000
$ws.Range("N121").Value = -36494

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4466.8096
$ws.Range("I132").Value = 1636.091
$ws.Range("J132").Value = 7580.6
$ws.Range("K132").Value = 4908.272999999999
$ws.Range("L132").Value = 22741.8
$ws.Range("M132").Value = -2378.272999999999
$ws.Range("N132").Value = -27801.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1660.6333
$ws.Range("I20").Value = 1090.1428
$ws.Range("J20").Value = 2159.8125
$ws.Range("K20").Value = 1090.1428
$ws.Range("L20").Value = 2159.8125
$ws.Range("M20").Value = -843.1428000000001
$ws.Range("N20").Value = -2653.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1199.4286
$ws.Range("I64").Value = 1448
$ws.Range("J64").Value = 1100
$ws.Range("K64").Value = 1448
$ws.Range("L64").Value = 1100
$ws.Range("M64").Value = -1223
$ws.Range("N64").Value = -1550

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 1199.4286
$ws.Range("I67").Value = 1448
$ws.Range("J67").Value = 1100
$ws.Range("K67").Value = 1448
$ws.Range("L67").Value = 1100
$ws.Range("M67").Value = -668
$ws.Range("N67").Value = -2660

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3940.7334
$ws.Range("I16").Value = 3988.875
$ws.Range("J16").Value = 3885.7144
$ws.Range("K16").Value = 3988.875
$ws.Range("L16").Value = 3885.7144
$ws.Range("M16").Value = -3701.875
$ws.Range("N16").Value = -4459.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7094360
$ws.Range("I31").Value = 1742.4117
$ws.Range("J31").Value = 25644282
$ws.Range("K31").Value = 1742.4117
$ws.Range("L31").Value = 25644282
$ws.Range("M31").Value = -1447.4117
$ws.Range("N31").Value = -25644872

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7094360
$ws.Range("I34").Value = 1742.4117
$ws.Range("J34").Value = 25644282
$ws.Range("K34").Value = 1742.4117
$ws.Range("L34").Value = 25644282
$ws.Range("M34").Value = -1540.4117
$ws.Range("N34").Value = -25644686

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 8500
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("N38").Value = -10754

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H46").Value = 8500
$ws.Range("J46").Value = 10000
$ws.Range("L46").Value = 10000
$ws.Range("N46").Value = -10422

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 3940.7334
$ws.Range("I113").Value = 3988.875
$ws.Range("J113").Value = 3885.7144
$ws.Range("K113").Value = 3988.875
$ws.Range("L113").Value = 3885.7144
$ws.Range("M113").Value = -1818.875
$ws.Range("N113").Value = -8225.714400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 138.94118
$ws.Range("I12").Value = 1.8
$ws.Range("J12").Value = 196.08333
$ws.Range("K12").Value = 5.4
$ws.Range("L12").Value = 588.24999
$ws.Range("M12").Value = 167.6
$ws.Range("N12").Value = -934.24999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2570.6875
$ws.Range("I70").Value = 1621.8334
$ws.Range("J70").Value = 3140
$ws.Range("K70").Value = 4865.5002
$ws.Range("L70").Value = 9420
$ws.Range("M70").Value = -4550.5002
$ws.Range("N70").Value = -10050

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 2570.6875
$ws.Range("I73").Value = 1621.8334
$ws.Range("J73").Value = 3140
$ws.Range("K73").Value = 4865.5002
$ws.Range("L73").Value = 9420
$ws.Range("M73").Value = -3773.5002
$ws.Range("N73").Value = -11604

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1155.2222
$ws.Range("I92").Value = 749.2
$ws.Range("J92").Value = 1662.75
$ws.Range("K92").Value = 2247.6
$ws.Range("L92").Value = 4988.25
$ws.Range("M92").Value = -999.6000000000004
$ws.Range("N92").Value = -7484.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 41667560
$ws.Range("I107").Value = 500000000
$ws.Range("J107").Value = 972.7273
$ws.Range("K107").Value = 1500000000
$ws.Range("L107").Value = 2918.1819
$ws.Range("M107").Value = -1499998080
$ws.Range("N107").Value = -6758.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4372.143
$ws.Range("I70").Value = 4150.8335
$ws.Range("J70").Value = 5700
$ws.Range("K70").Value = 4150.8335
$ws.Range("L70").Value = 5700
$ws.Range("M70").Value = -3880.8335
$ws.Range("N70").Value = -6240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4372.143
$ws.Range("I73").Value = 4150.8335
$ws.Range("J73").Value = 5700
$ws.Range("K73").Value = 4150.8335
$ws.Range("L73").Value = 5700
$ws.Range("M73").Value = -3214.8335
$ws.Range("N73").Value = -7572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6000.75
$ws.Range("I61").Value = 3889.889
$ws.Range("K61").Value = 3889.889
$ws.Range("M61").Value = -3687.889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6000.75
$ws.Range("I113").Value = 3889.889
$ws.Range("K113").Value = 3889.889
$ws.Range("M113").Value = -1719.889

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1041.5217
$ws.Range("I113").Value = 655
$ws.Range("J113").Value = 1148.8889
$ws.Range("K113").Value = 1965
$ws.Range("L113").Value = 3446.6667
$ws.Range("M113").Value = 205
$ws.Range("N113").Value = -7786.6667
